$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 111.90909
$ws.Range("I55").Value = 55.166668
$ws.Range("J55").Value = 180
$ws.Range("K55").Value = 55.166668
$ws.Range("L55").Value = 180
$ws.Range("M55").Value = 158.833332
$ws.Range("N55").Value = -608
$ws.Range("H64").Value = 2981.6667
$ws.Range("J64").Value = 2981.6667
$ws.Range("L64").Value = 2981.6667
$ws.Range("N64").Value = -3477.6667
$ws.Range("H67").Value = 2981.6667
$ws.Range("J67").Value = 2981.6667
$ws.Range("L67").Value = 2981.6667
$ws.Range("N67").Value = -4697.6667
$ws.Range("H76").Value = 38676.605
$ws.Range("J76").Value = 2747.3333
$ws.Range("L76").Value = 2747.3333
$ws.Range("N76").Value = -3377.3333
$ws.Range("H79").Value = 38676.605
$ws.Range("J79").Value = 2747.3333
$ws.Range("L79").Value = 2747.3333
$ws.Range("N79").Value = -4931.3333
$ws.Range("H99").Value = 340.91666
$ws.Range("I99").Value = 209.1
$ws.Range("K99").Value = 627.3
$ws.Range("M99").Value = 870.7

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 871.8570999999999
$ws.Range("I45").Value = 848.1667
$ws.Range("K45").Value = 848.1667
$ws.Range("M45").Value = -471.1667
$ws.Range("H61").Value = 1437.1351
$ws.Range("I61").Value = 1242.4667
$ws.Range("J61").Value = 2271.4285
$ws.Range("K61").Value = 1242.4667
$ws.Range("L61").Value = 2271.4285
$ws.Range("M61").Value = -1030.4667
$ws.Range("N61").Value = -2695.4285
$ws.Range("H122").Value = 4840.793
$ws.Range("I122").Value = 1530.6818
$ws.Range("J122").Value = 15244
$ws.Range("K122").Value = 4592.0454
$ws.Range("L122").Value = 45732
$ws.Range("M122").Value = -2142.0454
$ws.Range("N122").Value = -50632
$ws.Range("H132").Value = 7026.2085
$ws.Range("I132").Value = 8585.857
$ws.Range("J132").Value = 4842.7
$ws.Range("K132").Value = 25757.571
$ws.Range("L132").Value = 14528.1
$ws.Range("M132").Value = -23227.571
$ws.Range("N132").Value = -19588.1
$ws.Range("H136").Value = 1437.1351
$ws.Range("I136").Value = 1242.4667
$ws.Range("J136").Value = 2271.4285
$ws.Range("K136").Value = 3727.4001
$ws.Range("L136").Value = 6814.2855
$ws.Range("M136").Value = -1177.4001
$ws.Range("N136").Value = -11914.2855

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1962.6897
$ws.Range("I105").Value = 1873.8889
$ws.Range("J105").Value = 2108
$ws.Range("K105").Value = 1873.8889
$ws.Range("L105").Value = 2108
$ws.Range("M105").Value = -126.8888999999999
$ws.Range("N105").Value = -5602
$ws.Range("H134").Value = 74578
$ws.Range("I134").Value = 74578
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 223734
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -221199
$ws.Range("N134").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1961.6666
$ws.Range("I16").Value = 1020
$ws.Range("J16").Value = 3138.75
$ws.Range("K16").Value = 1020
$ws.Range("L16").Value = 3138.75
$ws.Range("M16").Value = -733
$ws.Range("N16").Value = -3712.75
$ws.Range("H31").Value = 1665.8169
$ws.Range("I31").Value = 2770.842
$ws.Range("K31").Value = 2770.842
$ws.Range("M31").Value = -2475.842
$ws.Range("H34").Value = 1665.8169
$ws.Range("I34").Value = 2770.842
$ws.Range("K34").Value = 2770.842
$ws.Range("M34").Value = -2568.842
$ws.Range("H58").Value = 15152429
$ws.Range("I58").Value = 969.5
$ws.Range("K58").Value = 969.5
$ws.Range("M58").Value = -766.5
$ws.Range("H113").Value = 1961.6666
$ws.Range("I113").Value = 1020
$ws.Range("J113").Value = 3138.75
$ws.Range("K113").Value = 1020
$ws.Range("L113").Value = 3138.75
$ws.Range("M113").Value = 1150
$ws.Range("N113").Value = -7478.75
$ws.Range("H132").Value = 2243.8293
$ws.Range("I132").Value = 1867.8214
$ws.Range("J132").Value = 3053.6924
$ws.Range("K132").Value = 5603.4642
$ws.Range("L132").Value = 9161.0772
$ws.Range("M132").Value = -3073.4642
$ws.Range("N132").Value = -14221.0772
$ws.Range("H134").Value = 1070.3636
$ws.Range("I134").Value = 1092.4
$ws.Range("J134").Value = 850
$ws.Range("K134").Value = 3277.2
$ws.Range("L134").Value = 2550
$ws.Range("M134").Value = -742.2000000000003
$ws.Range("N134").Value = -7620
$ws.Range("H136").Value = 15152429
$ws.Range("I136").Value = 969.5
$ws.Range("K136").Value = 2908.5
$ws.Range("M136").Value = -358.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1738
$ws.Range("I34").Value = 700
$ws.Range("J34").Value = 1997.5
$ws.Range("K34").Value = 2100
$ws.Range("L34").Value = 5992.5
$ws.Range("M34").Value = -2016
$ws.Range("N34").Value = -6160.5
$ws.Range("H68").Value = 1399.7963
$ws.Range("I68").Value = 1181.2572
$ws.Range("J68").Value = 1802.3684
$ws.Range("K68").Value = 3543.7716
$ws.Range("L68").Value = 5407.1052
$ws.Range("M68").Value = -2732.7716
$ws.Range("N68").Value = -7029.1052
$ws.Range("H70").Value = 4800
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 6666.6665
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 19999.9995
$ws.Range("M70").Value = -5685
$ws.Range("N70").Value = -20629.9995
$ws.Range("H71").Value = 1399.7963
$ws.Range("I71").Value = 1181.2572
$ws.Range("J71").Value = 1802.3684
$ws.Range("K71").Value = 10631.3148
$ws.Range("L71").Value = 16221.3156
$ws.Range("M71").Value = -6575.3148
$ws.Range("N71").Value = -24333.3156
$ws.Range("H73").Value = 4800
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 6666.6665
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 19999.9995
$ws.Range("M73").Value = -4908
$ws.Range("N73").Value = -22183.9995
$ws.Range("H75").Value = 708.6667
$ws.Range("I75").Value = 708.6667
$ws.Range("K75").Value = 2126.0001
$ws.Range("M75").Value = -1128.0001
$ws.Range("H78").Value = 708.6667
$ws.Range("I78").Value = 708.6667
$ws.Range("K78").Value = 6378.0003
$ws.Range("M78").Value = -1386.0003
$ws.Range("H86").Value = 5337.5
$ws.Range("J86").Value = 7000
$ws.Range("L86").Value = 21000
$ws.Range("N86").Value = -23372
$ws.Range("H89").Value = 5337.5
$ws.Range("J89").Value = 7000
$ws.Range("L89").Value = 63000
$ws.Range("N89").Value = -74856
$ws.Range("H107").Value = 618.381
$ws.Range("I107").Value = 415.58334
$ws.Range("J107").Value = 888.7778
$ws.Range("K107").Value = 1246.75002
$ws.Range("L107").Value = 2666.3334
$ws.Range("M107").Value = 673.2499800000001
$ws.Range("N107").Value = -6506.3334
$ws.Range("H129").Value = 6945592
$ws.Range("I129").Value = 517
$ws.Range("J129").Value = 11906360
$ws.Range("K129").Value = 1551
$ws.Range("L129").Value = 35719080
$ws.Range("M129").Value = 3449
$ws.Range("N129").Value = -35729080
$ws.Range("H131").Value = 6990225
$ws.Range("I131").Value = 8680
$ws.Range("J131").Value = 17462542
$ws.Range("K131").Value = 26040
$ws.Range("L131").Value = 52387626
$ws.Range("M131").Value = -21000
$ws.Range("N131").Value = -52397706

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1200.3846
$ws.Range("I97").Value = 990.5
$ws.Range("J97").Value = 1900
$ws.Range("K97").Value = 990.5
$ws.Range("L97").Value = 1900
$ws.Range("M97").Value = -494.5
$ws.Range("N97").Value = -2892
$ws.Range("H113").Value = 27781842
$ws.Range("I113").Value = 62506136
$ws.Range("J113").Value = 2408
$ws.Range("K113").Value = 62506136
$ws.Range("L113").Value = 2408
$ws.Range("M113").Value = -62503966
$ws.Range("N113").Value = -6748
$ws.Range("H122").Value = 43482464
$ws.Range("I122").Value = 71433550
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 214300650
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -214298200
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 86025.414
$ws.Range("I132").Value = 107616.42
$ws.Range("J132").Value = 3979.6
$ws.Range("K132").Value = 322849.26
$ws.Range("L132").Value = 11938.8
$ws.Range("M132").Value = -320319.26
$ws.Range("N132").Value = -16998.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1999.421
$ws.Range("I40").Value = 1963
$ws.Range("J40").Value = 2039.8889
$ws.Range("K40").Value = 1963
$ws.Range("L40").Value = 2039.8889
$ws.Range("M40").Value = -1827
$ws.Range("N40").Value = -2311.8889
$ws.Range("H122").Value = 1928.2
$ws.Range("I122").Value = 1934.8235
$ws.Range("J122").Value = 1921.9445
$ws.Range("K122").Value = 5804.470499999999
$ws.Range("L122").Value = 5765.833500000001
$ws.Range("M122").Value = -3354.470499999999
$ws.Range("N122").Value = -10665.8335
$ws.Range("H132").Value = 11397.25
$ws.Range("I132").Value = 15897.786
$ws.Range("J132").Value = 5096.5
$ws.Range("K132").Value = 47693.358
$ws.Range("L132").Value = 15289.5
$ws.Range("M132").Value = -45163.358
$ws.Range("N132").Value = -20349.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1375.9286
$ws.Range("I132").Value = 1231.2424
$ws.Range("K132").Value = 3693.7272
$ws.Range("M132").Value = -1163.7272
$ws.Range("H136").Value = 3832.0908
$ws.Range("I136").Value = 4394.125
$ws.Range("K136").Value = 13182.375
$ws.Range("M136").Value = -10632.375
